$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 7847.3076
$ws.Range("I9").Value = 11277.333
$ws.Range("J9").Value = 129.75
$ws.Range("K9").Value = 11277.333
$ws.Range("L9").Value = 129.75
$ws.Range("M9").Value = -11108.333
$ws.Range("N9").Value = -467.75
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 19
$ws.Range("H19").Value = 1285.6666
$ws.Range("J19").Value = 1303.625
$ws.Range("L19").Value = 1303.625
$ws.Range("N19").Value = -1653.625
# Row 40
$ws.Range("H40").Value = 3124.4
$ws.Range("I40").Value = 3042.1428
$ws.Range("J40").Value = 3316.3333
$ws.Range("K40").Value = 3042.1428
$ws.Range("L40").Value = 3316.3333
$ws.Range("M40").Value = -2867.1428
$ws.Range("N40").Value = -3666.3333
# Row 70
$ws.Range("H70").Value = 3991.9512
$ws.Range("I70").Value = 1398.3438
$ws.Range("K70").Value = 4195.0314
$ws.Range("M70").Value = -3925.0314
# Row 73
$ws.Range("H73").Value = 3991.9512
$ws.Range("I73").Value = 1398.3438
$ws.Range("K73").Value = 4195.0314
$ws.Range("M73").Value = -3259.0314
# Row 80
$ws.Range("H80").Value = 591.8125
$ws.Range("J80").Value = 698.4545000000001
$ws.Range("L80").Value = 2095.3635
$ws.Range("N80").Value = -4091.3635
# Row 83
$ws.Range("H83").Value = 591.8125
$ws.Range("J83").Value = 698.4545000000001
$ws.Range("L83").Value = 6286.0905
$ws.Range("N83").Value = -16270.0905
# Row 112
$ws.Range("H112").Value = 2173.6667
$ws.Range("J112").Value = 2144
$ws.Range("L112").Value = 6432
$ws.Range("N112").Value = -8648
# Row 137
$ws.Range("H137").Value = 1212.8605
$ws.Range("I137").Value = 944.9143
$ws.Range("K137").Value = 2834.7429
$ws.Range("M137").Value = -284.7429000000002
# Row 138
$ws.Range("H138").Value = 3638.31
$ws.Range("J138").Value = 3999
$ws.Range("L138").Value = 11997
$ws.Range("N138").Value = -22277

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1558.9672
$ws.Range("I32").Value = 1052.7119
$ws.Range("K32").Value = 1052.7119
$ws.Range("M32").Value = -765.7119
# Row 74
$ws.Range("H74").Value = 2535.577
$ws.Range("I74").Value = 2603
$ws.Range("K74").Value = 2603
$ws.Range("M74").Value = -1729
# Row 77
$ws.Range("H77").Value = 2535.577
$ws.Range("I77").Value = 2603
$ws.Range("K77").Value = 13015
$ws.Range("M77").Value = -8647
# Row 97
$ws.Range("H97").Value = 1157.9166
$ws.Range("I97").Value = 622.3333
$ws.Range("K97").Value = 622.3333
$ws.Range("M97").Value = -126.3333
# Row 122
$ws.Range("H122").Value = 2335.5757
$ws.Range("J122").Value = 3538.6924
$ws.Range("L122").Value = 10616.0772
$ws.Range("N122").Value = -15516.0772

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 7426.44
$ws.Range("I86").Value = 2807.7896
$ws.Range("J86").Value = 22052.166
$ws.Range("K86").Value = 2807.7896
$ws.Range("L86").Value = 22052.166
$ws.Range("M86").Value = -1684.7896
$ws.Range("N86").Value = -24298.166
# Row 89
$ws.Range("H89").Value = 7426.44
$ws.Range("I89").Value = 2807.7896
$ws.Range("J89").Value = 22052.166
$ws.Range("K89").Value = 14038.948
$ws.Range("L89").Value = 110260.83
$ws.Range("M89").Value = -8422.948
$ws.Range("N89").Value = -121492.83
# Row 94
$ws.Range("H94").Value = 6803.9688
$ws.Range("I94").Value = 4136.0415
$ws.Range("J94").Value = 14807.75
$ws.Range("K94").Value = 4136.0415
$ws.Range("L94").Value = 14807.75
$ws.Range("M94").Value = -3685.0415
$ws.Range("N94").Value = -15709.75
# Row 134
$ws.Range("H134").Value = 4994.5386
$ws.Range("I134").Value = 4216.243
$ws.Range("J134").Value = 6914.3335
$ws.Range("K134").Value = 12648.729
$ws.Range("L134").Value = 20743.0005
$ws.Range("M134").Value = -10113.729
$ws.Range("N134").Value = -25813.0005

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 31
$ws.Range("H31").Value = 2543
$ws.Range("I31").Value = 1938.5555
$ws.Range("J31").Value = 3223
$ws.Range("K31").Value = 1938.5555
$ws.Range("L31").Value = 3223
$ws.Range("M31").Value = -1643.5555
$ws.Range("N31").Value = -3813
# Row 34
$ws.Range("H34").Value = 2543
$ws.Range("I34").Value = 1938.5555
$ws.Range("J34").Value = 3223
$ws.Range("K34").Value = 1938.5555
$ws.Range("L34").Value = 3223
$ws.Range("M34").Value = -1736.5555
$ws.Range("N34").Value = -3627
# Row 62
$ws.Range("H62").Value = 30759.75
$ws.Range("J62").Value = 39013
$ws.Range("L62").Value = 39013
$ws.Range("N62").Value = -40261
# Row 65
$ws.Range("H65").Value = 30759.75
$ws.Range("J65").Value = 39013
$ws.Range("L65").Value = 195065
$ws.Range("N65").Value = -201305
# Row 132
$ws.Range("H132").Value = 3324.6943
$ws.Range("I132").Value = 1770.6923
$ws.Range("J132").Value = 7365.1
$ws.Range("K132").Value = 5312.0769
$ws.Range("L132").Value = 22095.3
$ws.Range("M132").Value = -2782.0769
$ws.Range("N132").Value = -27155.3

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 960
$ws.Range("I14").Value = 960
$ws.Range("K14").Value = 2880
$ws.Range("M14").Value = -2707
# Row 64
$ws.Range("H64").Value = 3490347.5
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
# Row 67
$ws.Range("H67").Value = 3490347.5
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
# Row 107
$ws.Range("H107").Value = 579.82355
$ws.Range("I107").Value = 198.875
$ws.Range("J107").Value = 918.44446
$ws.Range("K107").Value = 596.625
$ws.Range("L107").Value = 2755.33338
$ws.Range("M107").Value = 1323.375
$ws.Range("N107").Value = -6595.33338

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 12628.143
$ws.Range("I3").Value = 3449
$ws.Range("J3").Value = 16299.8
$ws.Range("K3").Value = 3449
$ws.Range("L3").Value = 16299.8
$ws.Range("M3").Value = -3333
$ws.Range("N3").Value = -16531.8
# Row 13
$ws.Range("H13").Value = 870
$ws.Range("I13").Value = 883.3333
$ws.Range("J13").Value = 850
$ws.Range("K13").Value = 883.3333
$ws.Range("L13").Value = 850
$ws.Range("M13").Value = -744.3333
$ws.Range("N13").Value = -1128
# Row 80
$ws.Range("H80").Value = 2949.625
$ws.Range("I80").Value = 2400
$ws.Range("K80").Value = 2400
$ws.Range("M80").Value = -1402
# Row 83
$ws.Range("H83").Value = 2949.625
$ws.Range("I83").Value = 2400
$ws.Range("K83").Value = 12000
$ws.Range("M83").Value = -7008
# Row 132
$ws.Range("H132").Value = 2541.96
$ws.Range("I132").Value = 2323.0476
$ws.Range("K132").Value = 6969.1428
$ws.Range("M132").Value = -4439.1428

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1254.1333
$ws.Range("I22").Value = 916.1111
$ws.Range("J22").Value = 1761.1666
$ws.Range("K22").Value = 916.1111
$ws.Range("L22").Value = 1761.1666
$ws.Range("M22").Value = -621.1111
$ws.Range("N22").Value = -2351.1666
# Row 27
$ws.Range("H27").Value = 1254.1333
$ws.Range("I27").Value = 916.1111
$ws.Range("J27").Value = 1761.1666
$ws.Range("K27").Value = 916.1111
$ws.Range("L27").Value = 1761.1666
$ws.Range("M27").Value = -809.1111
$ws.Range("N27").Value = -1975.1666
# Row 61
$ws.Range("H61").Value = 2648.625
$ws.Range("I61").Value = 2648.625
$ws.Range("K61").Value = 2648.625
$ws.Range("M61").Value = -2446.625
# Row 76
$ws.Range("H76").Value = 64999.5
$ws.Range("J76").Value = 64999.5
$ws.Range("L76").Value = 64999.5
$ws.Range("N76").Value = -65675.5
# Row 79
$ws.Range("H79").Value = 64999.5
$ws.Range("J79").Value = 64999.5
$ws.Range("L79").Value = 64999.5
$ws.Range("N79").Value = -67339.5
# Row 82
$ws.Range("H82").Value = 2378.8333
$ws.Range("I82").Value = 2126.6667
$ws.Range("K82").Value = 2126.6667
$ws.Range("M82").Value = -1765.6667
# Row 85
$ws.Range("H85").Value = 2378.8333
$ws.Range("I85").Value = 2126.6667
$ws.Range("K85").Value = 2126.6667
$ws.Range("M85").Value = -878.6667000000002
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 113
$ws.Range("H113").Value = 2648.625
$ws.Range("I113").Value = 2648.625
$ws.Range("K113").Value = 2648.625
$ws.Range("M113").Value = -478.625
# Row 132
$ws.Range("H132").Value = 2111.8647
$ws.Range("I132").Value = 1757.6111
$ws.Range("J132").Value = 2447.4736
$ws.Range("K132").Value = 5272.8333
$ws.Range("L132").Value = 7342.4208
$ws.Range("M132").Value = -2742.8333
$ws.Range("N132").Value = -12402.4208

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 6114
$ws.Range("J23").Value = 29999
$ws.Range("L23").Value = 29999
$ws.Range("N23").Value = -30457
# Row 81
$ws.Range("H81").Value = 4631.154
$ws.Range("I81").Value = 5355.4546
$ws.Range("J81").Value = 647.5
$ws.Range("K81").Value = 10710.9092
$ws.Range("L81").Value = 1295
$ws.Range("M81").Value = -9649.9092
$ws.Range("N81").Value = -3417
# Row 84
$ws.Range("H84").Value = 4631.154
$ws.Range("I84").Value = 5355.4546
$ws.Range("J84").Value = 647.5
$ws.Range("K84").Value = 53554.546
$ws.Range("L84").Value = 6475
$ws.Range("M84").Value = -48250.546
$ws.Range("N84").Value = -17083
